$d = $word.ActiveDocument

$replacements = @(
    @("631×5=", "949×3="),
    @("989×4=", "780×3="),
    @("300×9=", "294×4="),
    @("115×4=", "329×8="),
    @("678×6=", "857×6="),
    @("269×6=", "826×8="),
    @("574×8=", "618×9="),
    @("932×7=", "601×4="),
    @("944×7=", "800×4="),
    @("400×9=", "857×6="),
    @("946×7=", "972×9="),
    @("531×6=", "658×2="),
    @("461×2=", "742×3="),
    @("536×9=", "498×3="),
    @("718×5=", "513×3="),
    @("929×2=", "554×9="),
    @("388×8=", "399×7="),
    @("713×2=", "919×7="),
    @("522×4=", "565×2="),
    @("890×9=", "258×5="),
    @("362×2=", "273×8="),
    @("366×4=", "964×7="),
    @("637×5=", "817×8="),
    @("498×2=", "672×8="),
    @("800×2=", "457×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
